# Fix descriptives example data set: add a new "pvord_pooled" (pooled ordinal)
# row right before the existing "pvkat_pooled" row, so the sheet now covers
# all the different (pooled) variable types: pv, pvord and pvkat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "pvkat_pooled" row in column A so the new row is inserted in the
# right place even if the sheet layout shifts around.
$targetRow = 4
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value -eq "pvkat_pooled") {
        $targetRow = $r
        break
    }
}

# Insert a new row above the "pvkat_pooled" row. This shifts that row (and any
# rows below it) down by one, preserving its content.
$ws.Rows.Item($targetRow).Insert()

# Populate the newly inserted row with the pooled ordinal variable info.
$ws.Cells.Item($targetRow, 1).Value = "pvord_pooled"
$ws.Cells.Item($targetRow, 2).Value = "-"
$ws.Cells.Item($targetRow, 3).Value = ""
$ws.Cells.Item($targetRow, 4).Value = "pvord_1,pvord_2,pvord_3,pvord_4,pvord_5"
